# Weekly update: rotate the sliding window of "Fruta, Femacal de La Calera - Coco"
# price records (columns D, M, N, O, P, S) and append the new record that
# drops off the top of the window as a new row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row([int]$row, $d, $m, $n, $o, $p, $s) {
    $ws.Cells.Item($row, 4).Value = $d    # D - Fecha
    $ws.Cells.Item($row, 13).Value = $m   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $n   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $o   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $p   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $s   # S - Precio $/Kg
}

Set-Row 2  44307 30 22000 22000 22000 1100
Set-Row 3  44377 25 20000 20000 20000 1000
Set-Row 4  44305 20 22000 22000 22000 1100
Set-Row 5  44389 20 20000 20000 20000 1000
Set-Row 6  44382 24 20000 20000 20000 1000
Set-Row 7  44298 65 22000 22000 22000 1100
Set-Row 8  44445 45 20000 20000 20000 1000
Set-Row 9  44403 50 20000 20000 20000 1000
Set-Row 10 44413 45 20000 20000 20000 1000
Set-Row 11 44406 20 20000 20000 20000 1000
Set-Row 12 44385 36 20000 20000 20000 1000
Set-Row 13 44292 30 25000 25000 25000 1250
Set-Row 14 44300 45 22000 22000 22000 1100
Set-Row 15 44291 70 25000 25000 25000 1250
# Row 16 unchanged: D=44301 M=38 N=O=P=22000 S=1100
# Row 18 unchanged: D=44376 M=38 N=O=P=20000 S=1000
Set-Row 17 44294 25 25000 25000 25000 1250

# Append new row 19, a duplicate of the record that previously lived in row 8
# (Fecha 44400, Volumen 45, 20000/20000/20000, 1000 $/Kg) before it was
# replaced by the newer reading now stored in row 8.
$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44400
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100108
$ws.Cells.Item(19, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(19, 9).Value = 100108007
$ws.Cells.Item(19, 10).Value = "Coco"
$ws.Cells.Item(19, 11).Value = "Sin especificar"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 45
$ws.Cells.Item(19, 14).Value = 20000
$ws.Cells.Item(19, 15).Value = 20000
$ws.Cells.Item(19, 16).Value = 20000
$ws.Cells.Item(19, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(19, 18).Value = "Perú"
$ws.Cells.Item(19, 19).Value = 1000
$ws.Cells.Item(19, 20).Value = 20
